$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell N1: "Team", formatted like the other header cells
$ws.Range("N1").Value = "Team"
$ws.Range("A1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# Team column values for rows 2-121 (5-player blocks per team)
$ws.Range("N2:N6").Value = "GEN"
$ws.Range("N7:N11").Value = "RRQ"
$ws.Range("N12:N16").Value = "GEN"
$ws.Range("N17:N21").Value = "RRQ"
$ws.Range("N22:N26").Value = "T1"
$ws.Range("N27:N31").Value = "BLD"
$ws.Range("N32:N36").Value = "T1"
$ws.Range("N37:N41").Value = "BLD"
$ws.Range("N42:N46").Value = "ZETA"
$ws.Range("N47:N51").Value = "GE"
$ws.Range("N52:N56").Value = "ZETA"
$ws.Range("N57:N61").Value = "GE"
$ws.Range("N62:N66").Value = "BLD"
$ws.Range("N67:N71").Value = "GE"
$ws.Range("N72:N76").Value = "BLD"
$ws.Range("N77:N81").Value = "GE"
$ws.Range("N82:N86").Value = "T1"
$ws.Range("N87:N91").Value = "ZETA"
$ws.Range("N92:N96").Value = "T1"
$ws.Range("N97:N106").Value = "ZETA"
$ws.Range("N107:N111").Value = "GE"
$ws.Range("N112:N116").Value = "ZETA"
$ws.Range("N117:N121").Value = "GE"
